$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5353.684
$ws.Range("J32").Value = 4984.4443
$ws.Range("L32").Value = 4984.4443
$ws.Range("N32").Value = -5636.4443

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 919.76666
$ws.Range("I80").Value = 450.33334
$ws.Range("J80").Value = 1389.2
$ws.Range("K80").Value = 1351.00002
$ws.Range("L80").Value = 4167.6
$ws.Range("M80").Value = -353.0000199999999
$ws.Range("N80").Value = -6163.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 919.76666
$ws.Range("I83").Value = 450.33334
$ws.Range("J83").Value = 1389.2
$ws.Range("K83").Value = 4053.00006
$ws.Range("L83").Value = 12502.8
$ws.Range("M83").Value = 938.9999399999997
$ws.Range("N83").Value = -22486.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 50004028
$ws.Range("I132").Value = 58827740
$ws.Range("J132").Value = 2982.3333
$ws.Range("K132").Value = 176483220
$ws.Range("L132").Value = 8946.999899999999
$ws.Range("M132").Value = -176480690
$ws.Range("N132").Value = -14006.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 164975.1
$ws.Range("I137").Value = 357195.2
$ws.Range("J137").Value = 4791.6665
$ws.Range("K137").Value = 1071585.6
$ws.Range("L137").Value = 14374.9995
$ws.Range("M137").Value = -1069035.6
$ws.Range("N137").Value = -19474.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5499.7036
$ws.Range("I32").Value = 3369.9023
$ws.Range("K32").Value = 3369.9023
$ws.Range("M32").Value = -3082.9023

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2292.1177
$ws.Range("I61").Value = 2146.762
$ws.Range("K61").Value = 2146.762
$ws.Range("M61").Value = -1934.762

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4169881.8
$ws.Range("I102").Value = 5557821.5
$ws.Range("K102").Value = 5557821.5
$ws.Range("M102").Value = -5556199.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2299.457
$ws.Range("I132").Value = 2063.6667
$ws.Range("K132").Value = 6191.000100000001
$ws.Range("M132").Value = -3661.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2292.1177
$ws.Range("I136").Value = 2146.762
$ws.Range("K136").Value = 6440.286
$ws.Range("M136").Value = -3890.286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4332998
$ws.Range("I99").Value = 5956064.5
$ws.Range("K99").Value = 5956064.5
$ws.Range("M99").Value = -5954566.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 373.5
$ws.Range("I7").Value = 247.58333
$ws.Range("K7").Value = 247.58333
$ws.Range("M7").Value = -134.58333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 276.22223
$ws.Range("J22").Value = 328.33334
$ws.Range("L22").Value = 328.33334
$ws.Range("N22").Value = -1028.33334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15907.324
$ws.Range("I31").Value = 2271.3333
$ws.Range("J31").Value = 18424.738
$ws.Range("K31").Value = 2271.3333
$ws.Range("L31").Value = 18424.738
$ws.Range("M31").Value = -1976.3333
$ws.Range("N31").Value = -19014.738

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 15907.324
$ws.Range("I34").Value = 2271.3333
$ws.Range("J34").Value = 18424.738
$ws.Range("K34").Value = 2271.3333
$ws.Range("L34").Value = 18424.738
$ws.Range("M34").Value = -2069.3333
$ws.Range("N34").Value = -18828.738

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3605.6667
$ws.Range("I58").Value = 2996.1667
$ws.Range("K58").Value = 2996.1667
$ws.Range("M58").Value = -2793.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4462.4443
$ws.Range("I99").Value = 4440.5
$ws.Range("J99").Value = 4480
$ws.Range("K99").Value = 4440.5
$ws.Range("L99").Value = 4480
$ws.Range("M99").Value = -2942.5
$ws.Range("N99").Value = -7476

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4462.4443
$ws.Range("I126").Value = 4440.5
$ws.Range("J126").Value = 4480
$ws.Range("K126").Value = 13321.5
$ws.Range("L126").Value = 13440
$ws.Range("M126").Value = -10851.5
$ws.Range("N126").Value = -18380

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 39366.082
$ws.Range("I132").Value = 2359
$ws.Range("J132").Value = 446444
$ws.Range("K132").Value = 7077
$ws.Range("L132").Value = 1339332
$ws.Range("M132").Value = -4547
$ws.Range("N132").Value = -1344392

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3605.6667
$ws.Range("I136").Value = 2996.1667
$ws.Range("K136").Value = 8988.500100000001
$ws.Range("M136").Value = -6438.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 337.9355
$ws.Range("I2").Value = 48.555557
$ws.Range("K2").Value = 291.333342
$ws.Range("M2").Value = -178.333342

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 269.64706
$ws.Range("I23").Value = 20
$ws.Range("K23").Value = 60
$ws.Range("M23").Value = 175

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 6899.125
$ws.Range("I94").Value = 2174
$ws.Range("J94").Value = 8474.166999999999
$ws.Range("K94").Value = 6522
$ws.Range("L94").Value = 25422.501
$ws.Range("M94").Value = -5846
$ws.Range("N94").Value = -26774.501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1678.6428
$ws.Range("I107").Value = 1002
$ws.Range("J107").Value = 1863.1818
$ws.Range("K107").Value = 3006
$ws.Range("L107").Value = 5589.5454
$ws.Range("M107").Value = -1086
$ws.Range("N107").Value = -9429.545399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4139.222
$ws.Range("I132").Value = 3614.6365
$ws.Range("J132").Value = 4963.5713
$ws.Range("K132").Value = 10843.9095
$ws.Range("L132").Value = 14890.7139
$ws.Range("M132").Value = -8313.9095
$ws.Range("N132").Value = -19950.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I100").Value = 2991.5417
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2991.5417
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2450.5417
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7014.4546
$ws.Range("I132").Value = 6949.3125
$ws.Range("K132").Value = 20847.9375
$ws.Range("M132").Value = -18317.9375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 35606.324
$ws.Range("I136").Value = 44616.582
$ws.Range("K136").Value = 133849.746
$ws.Range("M136").Value = -131299.746

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 2999.5
$ws.Range("I8").Value = 3000
$ws.Range("K8").Value = 3000
$ws.Range("M8").Value = -2860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 40041948
$ws.Range("I132").Value = 47626164
$ws.Range("J132").Value = 224834.5
$ws.Range("K132").Value = 142878492
$ws.Range("L132").Value = 674503.5
$ws.Range("M132").Value = -142875962
$ws.Range("N132").Value = -679563.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I136").Value = 2971.375
$ws.Range("J136").Value = 5780
$ws.Range("K136").Value = 8914.125
$ws.Range("L136").Value = 17340
$ws.Range("M136").Value = -6364.125
$ws.Range("N136").Value = -22440
